# New .ttl from Google sheet has been generated.
#
# The regenerated sheet inserts a "dct:source(separator=",")" mapping column
# right after skos:definition@en in the header row (row 17), which pushes
# every following header one column to the right, AND drops the
# skos:broadMatch / skos:narrowMatch mapping columns entirely. Row 18's lone
# " " placeholder shifts along with it. Finally, the sheet no longer reaches
# out to columns X:Y, so those get cleared.
#
# NOTE: this runtime's Range.Insert/Range.Delete with a shift direction
# operate on the full column/row rather than just the targeted cell, so the
# header row is rebuilt here with plain per-cell value assignment instead
# (which only ever touches the single addressed cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: header / predicate-mapping row -------------------------------
$ws.Range("E17").Value = 'dct:source(separator=",")'
$ws.Range("F17").Value = 'skos:broader(separator=",")'
$ws.Range("G17").Value = 'skos:closeMatch(separator=",")'
$ws.Range("H17").Value = 'skos:exactMatch(separator=",")'
$ws.Range("I17").Value = 'skos:relatedMatch(separator=",")'
$ws.Range("J17").Value = 'owl:deprecated^^xsd:boolean'
$ws.Range("K17").Value = 'dct:isReplacedBy'
$ws.Range("L17").Value = 'skos:editorialNote@en'
$ws.Range("M17").Value = 'dct:creator(separator=",")'
$ws.Range("N17").Value = 'dct:contributor(separator=",")'
$ws.Range("O17").Value = ""

# --- Row 18: single vocabulary-term data row -------------------------------
# The lone " " placeholder value shifts one column to the right, from G18 to
# H18 (consistent with the column insert at E on row 17).
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = " "

# --- Drop now-empty trailing columns X:Y ------------------------------------
# The regenerated sheet no longer reaches out to X/Y; clearing them lets the
# used range / dimension shrink back down to column W.
$ws.Range("X1:Y18").Clear()
